$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 324; existing rows 324-343 shift down to 325-344.
$ws.Rows.Item(324).Insert()

# Populate the newly inserted row 324 with the new data record.
$ws.Cells.Item(324, 1).Value = 4
$ws.Cells.Item(324, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(324, 3).Value = "Los Lagos"
$ws.Cells.Item(324, 4).Value = 44585
$ws.Cells.Item(324, 5).Value = 10
$ws.Cells.Item(324, 6).Value = "Fruta"
$ws.Cells.Item(324, 7).Value = 100102
$ws.Cells.Item(324, 8).Value = "Cítricos"
$ws.Cells.Item(324, 9).Value = 100102005
$ws.Cells.Item(324, 10).Value = "Naranja"
$ws.Cells.Item(324, 11).Value = "Valencia"
$ws.Cells.Item(324, 12).Value = "Primera"
$ws.Cells.Item(324, 13).Value = 400
$ws.Cells.Item(324, 14).Value = 16000
$ws.Cells.Item(324, 15).Value = 16000
$ws.Cells.Item(324, 16).Value = 16000
$ws.Cells.Item(324, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(324, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(324, 19).Value = 1067
$ws.Cells.Item(324, 20).Value = 15
